$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor
$ws.Range("B3").Value = 0.9941926134706534
$ws.Range("C3").Value = 0.9944169140506466
$ws.Range("D3").Value = 0.9944757551700649

# Row 4: GradientBoostingRegressor
$ws.Range("B4").Value = 0.9956095007234073
$ws.Range("C4").Value = 0.9956093491478256
$ws.Range("D4").Value = 0.9956093390581503

# Row 5: AdaBoostRegressor
$ws.Range("B5").Value = 0.962151635987567
$ws.Range("C5").Value = 0.9655057756592421
$ws.Range("D5").Value = 0.9682829958002764
